$wb = $excel.ActiveWorkbook

# Sheet 1: LP1912
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 05:20:30"
$ws1.Range("A3").Value = "Total filas: 26"

$ws1.Cells.Item(6, 1).Value = "04:03:00"
$ws1.Cells.Item(6, 2).Value = "04:03"
$ws1.Cells.Item(6, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(6, 4).Value = 0
$ws1.Cells.Item(6, 5).Value = "LP1912"
$ws1.Cells.Item(7, 1).Value = "04:37:19"
$ws1.Cells.Item(7, 2).Value = "04:46"
$ws1.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(7, 4).Value = 9
$ws1.Cells.Item(7, 5).Value = "LP1912"
$ws1.Cells.Item(8, 1).Value = "04:52:25"
$ws1.Cells.Item(8, 2).Value = "04:53"
$ws1.Cells.Item(8, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(8, 4).Value = 1
$ws1.Cells.Item(8, 5).Value = "LP1912"
$ws1.Cells.Item(9, 1).Value = "04:52:25"
$ws1.Cells.Item(9, 2).Value = "05:16"
$ws1.Cells.Item(9, 3).Value = "17_ROMERO"
$ws1.Cells.Item(9, 4).Value = 24
$ws1.Cells.Item(9, 5).Value = "LP1912"
$ws1.Cells.Item(10, 1).Value = "05:20:30"
$ws1.Cells.Item(10, 2).Value = "05:20"
$ws1.Cells.Item(10, 3).Value = "17_ROMERO"
$ws1.Cells.Item(10, 4).Value = 0
$ws1.Cells.Item(10, 5).Value = "LP1912"
$ws1.Cells.Item(11, 1).Value = "04:52:25"
$ws1.Cells.Item(11, 2).Value = "05:22"
$ws1.Cells.Item(11, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(11, 4).Value = 30
$ws1.Cells.Item(11, 5).Value = "LP1912"
$ws1.Cells.Item(12, 1).Value = "05:20:30"
$ws1.Cells.Item(12, 2).Value = "05:26"
$ws1.Cells.Item(12, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(12, 4).Value = 6
$ws1.Cells.Item(12, 5).Value = "LP1912"
$ws1.Cells.Item(13, 1).Value = "05:20:30"
$ws1.Cells.Item(13, 2).Value = "05:34"
$ws1.Cells.Item(13, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(13, 4).Value = 14
$ws1.Cells.Item(13, 5).Value = "LP1912"
$ws1.Cells.Item(14, 1).Value = "04:03:00"
$ws1.Cells.Item(14, 2).Value = "05:35"
$ws1.Cells.Item(14, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(14, 4).Value = 92
$ws1.Cells.Item(14, 5).Value = "LP1912"
$ws1.Cells.Item(15, 1).Value = "04:03:00"
$ws1.Cells.Item(15, 2).Value = "05:41"
$ws1.Cells.Item(15, 3).Value = "14_ABASTO"
$ws1.Cells.Item(15, 4).Value = 98
$ws1.Cells.Item(15, 5).Value = "LP1912"
$ws1.Cells.Item(16, 1).Value = "05:20:30"
$ws1.Cells.Item(16, 2).Value = "05:46"
$ws1.Cells.Item(16, 3).Value = "15_ABASTO"
$ws1.Cells.Item(16, 4).Value = 26
$ws1.Cells.Item(16, 5).Value = "LP1912"
$ws1.Cells.Item(17, 1).Value = "05:20:30"
$ws1.Cells.Item(17, 2).Value = "05:54"
$ws1.Cells.Item(17, 3).Value = "10_OLMOS"
$ws1.Cells.Item(17, 4).Value = 34
$ws1.Cells.Item(17, 5).Value = "LP1912"
$ws1.Cells.Item(18, 1).Value = "05:20:30"
$ws1.Cells.Item(18, 2).Value = "06:04"
$ws1.Cells.Item(18, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(18, 4).Value = 44
$ws1.Cells.Item(18, 5).Value = "LP1912"
$ws1.Cells.Item(19, 1).Value = "05:20:30"
$ws1.Cells.Item(19, 2).Value = "06:11"
$ws1.Cells.Item(19, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(19, 4).Value = 51
$ws1.Cells.Item(19, 5).Value = "LP1912"
$ws1.Cells.Item(20, 1).Value = "05:20:30"
$ws1.Cells.Item(20, 2).Value = "06:14"
$ws1.Cells.Item(20, 3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(20, 4).Value = 54
$ws1.Cells.Item(20, 5).Value = "LP1912"
$ws1.Cells.Item(21, 1).Value = "05:20:30"
$ws1.Cells.Item(21, 2).Value = "06:21"
$ws1.Cells.Item(21, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(21, 4).Value = 61
$ws1.Cells.Item(21, 5).Value = "LP1912"
$ws1.Cells.Item(22, 1).Value = "05:20:30"
$ws1.Cells.Item(22, 2).Value = "06:27"
$ws1.Cells.Item(22, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(22, 4).Value = 67
$ws1.Cells.Item(22, 5).Value = "LP1912"
$ws1.Cells.Item(23, 1).Value = "05:20:30"
$ws1.Cells.Item(23, 2).Value = "06:29"
$ws1.Cells.Item(23, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(23, 4).Value = 69
$ws1.Cells.Item(23, 5).Value = "LP1912"
$ws1.Cells.Item(24, 1).Value = "05:20:30"
$ws1.Cells.Item(24, 2).Value = "06:31"
$ws1.Cells.Item(24, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(24, 4).Value = 71
$ws1.Cells.Item(24, 5).Value = "LP1912"
$ws1.Cells.Item(25, 1).Value = "05:20:30"
$ws1.Cells.Item(25, 2).Value = "06:44"
$ws1.Cells.Item(25, 3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(25, 4).Value = 84
$ws1.Cells.Item(25, 5).Value = "LP1912"
$ws1.Cells.Item(26, 1).Value = "05:20:30"
$ws1.Cells.Item(26, 2).Value = "06:46"
$ws1.Cells.Item(26, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(26, 4).Value = 86
$ws1.Cells.Item(26, 5).Value = "LP1912"
$ws1.Cells.Item(27, 1).Value = "05:20:30"
$ws1.Cells.Item(27, 2).Value = "06:59"
$ws1.Cells.Item(27, 3).Value = "14_ABASTO"
$ws1.Cells.Item(27, 4).Value = 99
$ws1.Cells.Item(27, 5).Value = "LP1912"
$ws1.Cells.Item(28, 1).Value = "05:20:30"
$ws1.Cells.Item(28, 2).Value = "07:05"
$ws1.Cells.Item(28, 3).Value = "15_ABASTO"
$ws1.Cells.Item(28, 4).Value = 105
$ws1.Cells.Item(28, 5).Value = "LP1912"
$ws1.Cells.Item(29, 1).Value = "05:20:30"
$ws1.Cells.Item(29, 2).Value = "07:07"
$ws1.Cells.Item(29, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(29, 4).Value = 107
$ws1.Cells.Item(29, 5).Value = "LP1912"
$ws1.Cells.Item(30, 1).Value = "05:20:30"
$ws1.Cells.Item(30, 2).Value = "07:11"
$ws1.Cells.Item(30, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(30, 4).Value = 111
$ws1.Cells.Item(30, 5).Value = "LP1912"
$ws1.Cells.Item(31, 1).Value = "05:20:30"
$ws1.Cells.Item(31, 2).Value = "07:15"
$ws1.Cells.Item(31, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(31, 4).Value = 115
$ws1.Cells.Item(31, 5).Value = "LP1912"

# Sheet 2: LP1912-215
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 05:20:30"
$ws2.Range("A3").Value = "Total filas: 6"

$ws2.Cells.Item(6, 1).Value = "04:37:19"
$ws2.Cells.Item(6, 2).Value = "04:46"
$ws2.Cells.Item(6, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(6, 4).Value = 9
$ws2.Cells.Item(6, 5).Value = "LP1912"
$ws2.Cells.Item(7, 1).Value = "05:20:30"
$ws2.Cells.Item(7, 2).Value = "05:34"
$ws2.Cells.Item(7, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(7, 4).Value = 14
$ws2.Cells.Item(7, 5).Value = "LP1912"
$ws2.Cells.Item(8, 1).Value = "04:03:00"
$ws2.Cells.Item(8, 2).Value = "05:35"
$ws2.Cells.Item(8, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(8, 4).Value = 92
$ws2.Cells.Item(8, 5).Value = "LP1912"
$ws2.Cells.Item(9, 1).Value = "05:20:30"
$ws2.Cells.Item(9, 2).Value = "06:11"
$ws2.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(9, 4).Value = 51
$ws2.Cells.Item(9, 5).Value = "LP1912"
$ws2.Cells.Item(10, 1).Value = "05:20:30"
$ws2.Cells.Item(10, 2).Value = "06:46"
$ws2.Cells.Item(10, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(10, 4).Value = 86
$ws2.Cells.Item(10, 5).Value = "LP1912"
$ws2.Cells.Item(11, 1).Value = "05:20:30"
$ws2.Cells.Item(11, 2).Value = "07:11"
$ws2.Cells.Item(11, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(11, 4).Value = 111
$ws2.Cells.Item(11, 5).Value = "LP1912"

# Sheet 3: 6203-6173
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A1").Value = "LÍNEA 141 - 6203-6173 - 27/01/2026"
$ws3.Range("A2").Value = "Última actualización: 05:20:30"

$ws3.Cells.Item(6, 1).Value = "05:20:30"
$ws3.Cells.Item(6, 2).Value = "05:44"
$ws3.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6, 4).Value = 24
$ws3.Cells.Item(6, 5).Value = "L6173"
$ws3.Cells.Item(7, 1).Value = "04:52:25"
$ws3.Cells.Item(7, 2).Value = "06:09"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 77
$ws3.Cells.Item(7, 5).Value = "L6173"
$ws3.Cells.Item(8, 1).Value = "05:20:30"
$ws3.Cells.Item(8, 2).Value = "06:10"
$ws3.Cells.Item(8, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(8, 4).Value = 50
$ws3.Cells.Item(8, 5).Value = "L6173"
$ws3.Cells.Item(9, 1).Value = "05:20:30"
$ws3.Cells.Item(9, 2).Value = "06:33"
$ws3.Cells.Item(9, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(9, 4).Value = 73
$ws3.Cells.Item(9, 5).Value = "L6203"
$ws3.Cells.Item(10, 1).Value = "05:20:30"
$ws3.Cells.Item(10, 2).Value = "07:00"
$ws3.Cells.Item(10, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(10, 4).Value = 100
$ws3.Cells.Item(10, 5).Value = "L6173"

